# Balance.xlsx update:
# The Xcg estimation-method comparison rows for SFORZA and TORENBEEK_1982
# swapped places (with their values staying attached to their own label)
# on the FUSELAGE sheet (rows 23/24) and on the WING sheet
# (rows 23/24 and rows 27/28).

$wb = $excel.ActiveWorkbook

# --- FUSELAGE sheet ---
$wsFuselage = $wb.Worksheets.Item("FUSELAGE")

$fusA23 = $wsFuselage.Range("A23").Value2
$fusC23 = $wsFuselage.Range("C23").Value2
$fusA24 = $wsFuselage.Range("A24").Value2
$fusC24 = $wsFuselage.Range("C24").Value2

$wsFuselage.Range("A23").Value2 = $fusA24
$wsFuselage.Range("C23").Value2 = $fusC24
$wsFuselage.Range("A24").Value2 = $fusA23
$wsFuselage.Range("C24").Value2 = $fusC23

# --- WING sheet ---
$wsWing = $wb.Worksheets.Item("WING")

$wingA23 = $wsWing.Range("A23").Value2
$wingC23 = $wsWing.Range("C23").Value2
$wingA24 = $wsWing.Range("A24").Value2
$wingC24 = $wsWing.Range("C24").Value2

$wsWing.Range("A23").Value2 = $wingA24
$wsWing.Range("C23").Value2 = $wingC24
$wsWing.Range("A24").Value2 = $wingA23
$wsWing.Range("C24").Value2 = $wingC23

$wingA27 = $wsWing.Range("A27").Value2
$wingC27 = $wsWing.Range("C27").Value2
$wingA28 = $wsWing.Range("A28").Value2
$wingC28 = $wsWing.Range("C28").Value2

$wsWing.Range("A27").Value2 = $wingA28
$wsWing.Range("C27").Value2 = $wingC28
$wsWing.Range("A28").Value2 = $wingA27
$wsWing.Range("C28").Value2 = $wingC27
